$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.077097194815224
$ws.Cells.Item(2, 4).Value = 1.086920381594432
$ws.Cells.Item(2, 5).Value = 1.081559366048157
$ws.Cells.Item(2, 6).Value = 1.094771762609407
$ws.Cells.Item(2, 9).Value = 1.02359499962809
$ws.Cells.Item(2, 10).Value = 1.081993357721316
$ws.Cells.Item(2, 11).Value = 1.089576980966636
$ws.Cells.Item(2, 12).Value = 1.084229872939025
$ws.Cells.Item(2, 13).Value = 1.097408259026155
$ws.Cells.Item(2, 14).Value = 1.030231680001454

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.07962987264386
$ws.Cells.Item(3, 4).Value = 1.089394839400523
$ws.Cells.Item(3, 5).Value = 1.083873224589529
$ws.Cells.Item(3, 6).Value = 1.09728640589567
$ws.Cells.Item(3, 9).Value = 1.023504579208684
$ws.Cells.Item(3, 10).Value = 1.08417879815469
$ws.Cells.Item(3, 11).Value = 1.091866945882649
$ws.Cells.Item(3, 12).Value = 1.086358584309259
$ws.Cells.Item(3, 13).Value = 1.099739819178849
$ws.Cells.Item(3, 14).Value = 1.031019165620977

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.081261425383543
$ws.Cells.Item(4, 4).Value = 1.090989084607024
$ws.Cells.Item(4, 5).Value = 1.085363605105732
$ws.Cells.Item(4, 6).Value = 1.09890668132633
$ws.Cells.Item(4, 9).Value = 1.023443335729026
$ws.Cells.Item(4, 10).Value = 1.085585655038973
$ws.Cells.Item(4, 11).Value = 1.093341504606984
$ws.Cells.Item(4, 12).Value = 1.087728826230066
$ws.Cells.Item(4, 13).Value = 1.101241321932016
$ws.Cells.Item(4, 14).Value = 1.031524763764052

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.081945641130516
$ws.Cells.Item(5, 4).Value = 1.091657701279923
$ws.Cells.Item(5, 5).Value = 1.085988567002672
$ws.Cells.Item(5, 6).Value = 1.099586248201469
$ws.Cells.Item(5, 9).Value = 1.023416932628352
$ws.Cells.Item(5, 10).Value = 1.086175400095348
$ws.Cells.Item(5, 11).Value = 1.093959729160351
$ws.Cells.Item(5, 12).Value = 1.08830319924295
$ws.Cells.Item(5, 13).Value = 1.101870881242414
$ws.Cells.Item(5, 14).Value = 1.031736384892076

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.082060426328191
$ws.Cells.Item(6, 4).Value = 1.091769872243304
$ws.Cells.Item(6, 5).Value = 1.086093408635332
$ws.Cells.Item(6, 6).Value = 1.099700258182677
$ws.Cells.Item(6, 9).Value = 1.023412460900047
$ws.Cells.Item(6, 10).Value = 1.086274322569442
$ws.Cells.Item(6, 11).Value = 1.094063434544672
$ws.Cells.Item(6, 12).Value = 1.088399541871634
$ws.Cells.Item(6, 13).Value = 1.101976490236201
$ws.Cells.Item(6, 14).Value = 1.031771862775023

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.081270574479414
$ws.Cells.Item(7, 4).Value = 1.090998024936305
$ws.Cells.Item(7, 5).Value = 1.085371962082602
$ws.Cells.Item(7, 6).Value = 1.098915767944191
$ws.Cells.Item(7, 9).Value = 1.02344298551034
$ws.Cells.Item(7, 10).Value = 1.085593541854633
$ws.Cells.Item(7, 11).Value = 1.093349771895775
$ws.Cells.Item(7, 12).Value = 1.087736507563336
$ws.Cells.Item(7, 13).Value = 1.101249740646279
$ws.Cells.Item(7, 14).Value = 1.031527595099063

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.077954657831719
$ws.Cells.Item(8, 4).Value = 1.087758091969406
$ws.Cells.Item(8, 5).Value = 1.082342788236642
$ws.Cells.Item(8, 6).Value = 1.095623048461126
$ws.Cells.Item(8, 9).Value = 1.023565007130014
$ws.Cells.Item(8, 10).Value = 1.082733468584535
$ws.Cells.Item(8, 11).Value = 1.090352403287819
$ws.Cells.Item(8, 12).Value = 1.084950792022394
$ws.Cells.Item(8, 13).Value = 1.098197732932352
$ws.Cells.Item(8, 14).Value = 1.030498643121221

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.072053773455259
$ws.Cells.Item(9, 4).Value = 1.081993998353002
$ws.Cells.Item(9, 5).Value = 1.076950624113912
$ws.Cells.Item(9, 6).Value = 1.089766126209512
$ws.Cells.Item(9, 9).Value = 1.023759139049229
$ws.Cells.Item(9, 10).Value = 1.077636040674979
$ws.Cells.Item(9, 11).Value = 1.085013493287195
$ws.Cells.Item(9, 12).Value = 1.079985158096746
$ws.Cells.Item(9, 13).Value = 1.092762741635001
$ws.Cells.Item(9, 14).Value = 1.028654497441378

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.068077840655972
$ws.Cells.Item(10, 4).Value = 1.078111355415872
$ws.Cells.Item(10, 5).Value = 1.073316488630716
$ws.Cells.Item(10, 6).Value = 1.085821691539815
$ws.Cells.Item(10, 9).Value = 1.023874603872166
$ws.Cells.Item(10, 10).Value = 1.074196275556257
$ws.Cells.Item(10, 11).Value = 1.081412947653394
$ws.Cells.Item(10, 12).Value = 1.076633851925781
$ws.Cells.Item(10, 13).Value = 1.089098243217101
$ws.Cells.Item(10, 14).Value = 1.027403246845007

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.066345544997928
$ws.Cells.Item(11, 4).Value = 1.076419985771268
$ws.Cells.Item(11, 5).Value = 1.071732899355283
$ws.Cells.Item(11, 6).Value = 1.084103580404519
$ws.Cells.Item(11, 9).Value = 1.023921309042974
$ws.Cells.Item(11, 10).Value = 1.072696360970624
$ws.Cells.Item(11, 11).Value = 1.079843441168943
$ws.Cells.Item(11, 12).Value = 1.075172400775276
$ws.Cells.Item(11, 13).Value = 1.087501061105503
$ws.Cells.Item(11, 14).Value = 1.026856044875795

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.065700425458806
$ws.Cells.Item(12, 4).Value = 1.075790150631391
$ws.Cells.Item(12, 5).Value = 1.071143128613405
$ws.Cells.Item(12, 6).Value = 1.083463813866095
$ws.Cells.Item(12, 9).Value = 1.023938164268118
$ws.Cells.Item(12, 10).Value = 1.072137598405157
$ws.Cells.Item(12, 11).Value = 1.079258831294015
$ws.Cells.Item(12, 12).Value = 1.074627950888349
$ws.Cells.Item(12, 13).Value = 1.086906173083127
$ws.Cells.Item(12, 14).Value = 1.026651958764889

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.065838882266833
$ws.Cells.Item(13, 4).Value = 1.075925325096326
$ws.Cells.Item(13, 5).Value = 1.071269707712809
$ws.Cells.Item(13, 6).Value = 1.083601118588456
$ws.Cells.Item(13, 9).Value = 1.023934571058996
$ws.Cells.Item(13, 10).Value = 1.072257529400073
$ws.Cells.Item(13, 11).Value = 1.079384306554306
$ws.Cells.Item(13, 12).Value = 1.074744810578688
$ws.Cells.Item(13, 13).Value = 1.087033852969395
$ws.Cells.Item(13, 14).Value = 1.026695773845095

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.06629225359982
$ws.Cells.Item(14, 4).Value = 1.076367956083231
$ws.Cells.Item(14, 5).Value = 1.071684180769943
$ws.Cells.Item(14, 6).Value = 1.084050729755418
$ws.Cells.Item(14, 9).Value = 1.023922712353277
$ws.Cells.Item(14, 10).Value = 1.072650206999685
$ws.Cells.Item(14, 11).Value = 1.079795150616917
$ws.Cells.Item(14, 12).Value = 1.07512742936232
$ws.Cells.Item(14, 13).Value = 1.087451920925396
$ws.Cells.Item(14, 14).Value = 1.026839192118332

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.066571367774509
$ws.Cells.Item(15, 4).Value = 1.076640463816594
$ws.Cells.Item(15, 5).Value = 1.071939343556082
$ws.Cells.Item(15, 6).Value = 1.084327538398557
$ws.Cells.Item(15, 9).Value = 1.023915340507297
$ws.Cells.Item(15, 10).Value = 1.072891931335021
$ws.Cells.Item(15, 11).Value = 1.080048068147763
$ws.Cells.Item(15, 12).Value = 1.075362959575611
$ws.Cells.Item(15, 13).Value = 1.087709289556107
$ws.Cells.Item(15, 14).Value = 1.026927446157972

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.068192576942072
$ws.Cells.Item(16, 4).Value = 1.078223387008178
$ws.Cells.Item(16, 5).Value = 1.073421371194162
$ws.Cells.Item(16, 6).Value = 1.085935498133599
$ws.Cells.Item(16, 9).Value = 1.023871434968932
$ws.Cells.Item(16, 10).Value = 1.074295594700403
$ws.Cells.Item(16, 11).Value = 1.081516885796698
$ws.Cells.Item(16, 12).Value = 1.076730621886965
$ws.Cells.Item(16, 13).Value = 1.089204018390464
$ws.Cells.Item(16, 14).Value = 1.027439447380139

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.069206613332419
$ws.Cells.Item(17, 4).Value = 1.07921355177558
$ws.Cells.Item(17, 5).Value = 1.074348296271926
$ws.Cells.Item(17, 6).Value = 1.086941371166386
$ws.Cells.Item(17, 9).Value = 1.023843013862359
$ws.Cells.Item(17, 10).Value = 1.07517323382403
$ws.Cells.Item(17, 11).Value = 1.082435400445559
$ws.Cells.Item(17, 12).Value = 1.077585722482316
$ws.Cells.Item(17, 13).Value = 1.090138790368431
$ws.Cells.Item(17, 14).Value = 1.027759152391226

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.069797056351687
$ws.Cells.Item(18, 4).Value = 1.0797901219054
$ws.Cells.Item(18, 5).Value = 1.074887996393725
$ws.Cells.Item(18, 6).Value = 1.087527105126903
$ws.Cells.Item(18, 9).Value = 1.023826118601784
$ws.Cells.Item(18, 10).Value = 1.075684139027546
$ws.Cells.Item(18, 11).Value = 1.082970150395692
$ws.Cells.Item(18, 12).Value = 1.078083496684618
$ws.Cells.Item(18, 13).Value = 1.090683024565266
$ws.Cells.Item(18, 14).Value = 1.027945111006379

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.069998209770692
$ws.Cells.Item(19, 4).Value = 1.079986553595488
$ws.Cells.Item(19, 5).Value = 1.075071859074925
$ws.Cells.Item(19, 6).Value = 1.087726661713932
$ws.Cells.Item(19, 9).Value = 1.023820303822167
$ws.Cells.Item(19, 10).Value = 1.075858175336502
$ws.Cells.Item(19, 11).Value = 1.083152317668115
$ws.Cells.Item(19, 12).Value = 1.078253058215992
$ws.Cells.Item(19, 13).Value = 1.090868426015081
$ws.Cells.Item(19, 14).Value = 1.028008430481286

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.069097923520204
$ws.Cells.Item(20, 4).Value = 1.079107417864153
$ws.Cells.Item(20, 5).Value = 1.074248945633899
$ws.Cells.Item(20, 6).Value = 1.086833551759724
$ws.Cells.Item(20, 9).Value = 1.023846096030143
$ws.Cells.Item(20, 10).Value = 1.075079175982128
$ws.Cells.Item(20, 11).Value = 1.082336956776091
$ws.Cells.Item(20, 12).Value = 1.077494081213012
$ws.Cells.Item(20, 13).Value = 1.090038602278289
$ws.Cells.Item(20, 14).Value = 1.027724904984341

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.06615879358326
$ws.Cells.Item(21, 4).Value = 1.076237656525637
$ws.Cells.Item(21, 5).Value = 1.071562172163622
$ws.Cells.Item(21, 6).Value = 1.083918374670249
$ws.Cells.Item(21, 9).Value = 1.02392621804881
$ws.Cells.Item(21, 10).Value = 1.072534618599807
$ws.Cells.Item(21, 11).Value = 1.079674212600262
$ws.Cells.Item(21, 12).Value = 1.075014802312157
$ws.Cells.Item(21, 13).Value = 1.08732885560738
$ws.Cells.Item(21, 14).Value = 1.026796982106494

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.064301165077549
$ws.Cells.Item(22, 4).Value = 1.074424122510582
$ws.Cells.Item(22, 5).Value = 1.069863865158007
$ws.Cells.Item(22, 6).Value = 1.082076294799277
$ws.Cells.Item(22, 9).Value = 1.023973741382756
$ws.Cells.Item(22, 10).Value = 1.07092531059579
$ws.Cells.Item(22, 11).Value = 1.077990607768552
$ws.Cells.Item(22, 12).Value = 1.073446686672098
$ws.Cells.Item(22, 13).Value = 1.08561570860124
$ws.Cells.Item(22, 14).Value = 1.026208742761087

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.065286867658085
$ws.Cells.Item(23, 4).Value = 1.075386403314035
$ws.Cells.Item(23, 5).Value = 1.070765044081621
$ws.Cells.Item(23, 6).Value = 1.083053707678578
$ws.Cells.Item(23, 9).Value = 1.023948818288664
$ws.Cells.Item(23, 10).Value = 1.071779348693187
$ws.Cells.Item(23, 11).Value = 1.078884031521249
$ws.Cells.Item(23, 12).Value = 1.074278873237088
$ws.Cells.Item(23, 13).Value = 1.086524792474834
$ws.Cells.Item(23, 14).Value = 1.026521042759848

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.069147038939072
$ws.Cells.Item(24, 4).Value = 1.079155378223981
$ws.Cells.Item(24, 5).Value = 1.07429384087325
$ws.Cells.Item(24, 6).Value = 1.086882273717744
$ws.Cells.Item(24, 9).Value = 1.023844704313072
$ws.Cells.Item(24, 10).Value = 1.075121679780783
$ws.Cells.Item(24, 11).Value = 1.082381442332873
$ws.Cells.Item(24, 12).Value = 1.077535493019211
$ws.Cells.Item(24, 13).Value = 1.090083876058649
$ws.Cells.Item(24, 14).Value = 1.027740381518993

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.073586464254422
$ws.Cells.Item(25, 4).Value = 1.083490967031524
$ws.Cells.Item(25, 5).Value = 1.078351356378972
$ws.Cells.Item(25, 6).Value = 1.091287074335799
$ws.Cells.Item(25, 9).Value = 1.02371141742819
$ws.Cells.Item(25, 10).Value = 1.078960954643524
$ws.Cells.Item(25, 11).Value = 1.086400792518705
$ws.Cells.Item(25, 12).Value = 1.081275901106545
$ws.Cells.Item(25, 13).Value = 1.094174859860201
$ws.Cells.Item(25, 14).Value = 1.029135026410762
